$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update "Förändrad" (C column) date serial from 45221 (2023-10-22) to 45224 (2023-10-25)
# for rows 2 through 7, preserving the existing date formatting/style of the cells.
foreach ($row in 2..7) {
    $ws.Cells.Item($row, 3).Value = 45224
}
